$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (column C) and p-values (column D)
$ws.Range("C2").Value = -1.245855044786252
$ws.Range("D2").Value = 0.2213345434467053

$ws.Range("C3").Value = -0.1118515859381203
$ws.Range("D3").Value = 0.9115988256749008

$ws.Range("C4").Value = 2.069143425949857
$ws.Range("D4").Value = 0.04620214098780884

$ws.Range("C5").Value = 0.4078047655749159
$ws.Range("D5").Value = 0.6859736425470171

$ws.Range("C6").Value = 0.9819764891855218
$ws.Range("D6").Value = 0.3330478072657197

$ws.Range("C7").Value = 4.407501941211613
$ws.Range("D7").Value = 0.00009929269700292842

$ws.Range("C8").Value = 1.577163185876448
$ws.Range("D8").Value = 0.1240181071225848

$ws.Range("C9").Value = 2.305394427256292
$ws.Range("D9").Value = 0.02737482591935736

$ws.Range("C10").Value = 0.8371709646527012
$ws.Range("D10").Value = 0.4083426236814005

$ws.Range("C11").Value = -1.982418149785275
$ws.Range("D11").Value = 0.05556059809184344

$wb.Save()
